$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.5022738159769057
$ws.Range("J2").Value = 0.5022738159769057
$ws.Range("P2").Value = 0.9810128591839974
$ws.Range("S2").Value = 0.4927370723047613
$ws.Range("T2").Value = 0.4927370723047612

# Row 3
$ws.Range("I3").Value = 0.5022738159769057
$ws.Range("J3").Value = 0.5022738159769057
$ws.Range("S3").Value = 0.009536743672144455
$ws.Range("T3").Value = 0.009536743672144455

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.63173
$ws.Range("H4").Value = 1.89519
$ws.Range("I4").Value = 0.4977261840230943
$ws.Range("J4").Value = 0.4977261840230943
$ws.Range("P4").Value = 0.9810128591839974
$ws.Range("Q4").Value = 0.006506187270000001
$ws.Range("R4").Value = 0.05855568543
$ws.Range("S4").Value = 0.4882757868792362
$ws.Range("T4").Value = 0.4882757868792362

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.63173
$ws.Range("H5").Value = 1.89519
$ws.Range("I5").Value = 0.4977261840230943
$ws.Range("J5").Value = 0.4977261840230943
$ws.Range("Q5").Value = 0.0001259248466666667
$ws.Range("R5").Value = 0.00113332362
$ws.Range("S5").Value = 0.009450397143858084
$ws.Range("T5").Value = 0.009450397143858084
